# Apply updated Newey-West standard-error-adjusted correlation matrix values
# to the corrM correlation table (WorkingFolder/Tables/corrM.xlsx -> Sheet1).
# Every off-diagonal cell in the symmetric 14x14 matrix is rewritten; the
# diagonal (value 1) and all other sheet content are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = -0.6308995147500133
$ws.Range("D2").Value = 0.0347600473387599
$ws.Range("E2").Value = 0.03431866192840193
$ws.Range("G2").Value = -0.02527535427880465
$ws.Range("H2").Value = -0.3848572087141396
$ws.Range("I2").Value = 0.07613113137254446
$ws.Range("J2").Value = 0.09692859617640064
$ws.Range("K2").Value = 0.1457023439347263
$ws.Range("L2").Value = 0.06131781591624236
$ws.Range("M2").Value = 0.06612057560974671
$ws.Range("N2").Value = -0.3919766247436938
$ws.Range("O2").Value = 0.1208259260153849

# Row 3
$ws.Range("B3").Value = -0.6308995147500133
$ws.Range("D3").Value = 0.07440478395551647
$ws.Range("E3").Value = -0.06986646035450562
$ws.Range("G3").Value = -0.06051795811232875
$ws.Range("H3").Value = 0.1143498427673656
$ws.Range("I3").Value = 0.09527665589281313
$ws.Range("J3").Value = -0.01151755191354742
$ws.Range("K3").Value = -0.02914244125909981
$ws.Range("L3").Value = 0.04876962339280161
$ws.Range("M3").Value = 0.01742559811962423
$ws.Range("N3").Value = 0.06934592894886464
$ws.Range("O3").Value = 0.06306186227928907

# Row 4
$ws.Range("B4").Value = 0.0347600473387599
$ws.Range("C4").Value = 0.07440478395551647
$ws.Range("E4").Value = 0.3456618958283406
$ws.Range("G4").Value = 0.3269328383725711
$ws.Range("H4").Value = 0.3958692142618715
$ws.Range("I4").Value = 0.03844707986648056
$ws.Range("J4").Value = 0.8041894500608463
$ws.Range("K4").Value = 0.1418704476126033
$ws.Range("L4").Value = -0.01463024890903251
$ws.Range("M4").Value = 0.05679644052597838
$ws.Range("N4").Value = 0.3661331865400339
$ws.Range("O4").Value = 0.04245051279411372

# Row 5
$ws.Range("B5").Value = 0.03431866192840193
$ws.Range("C5").Value = -0.06986646035450562
$ws.Range("D5").Value = 0.3456618958283406
$ws.Range("G5").Value = 0.952956949937701
$ws.Range("H5").Value = 0.3616315137833193
$ws.Range("I5").Value = 0.3244402953268535
$ws.Range("J5").Value = 0.6315695168996384
$ws.Range("K5").Value = 0.06432798670309411
$ws.Range("L5").Value = -0.09781173102871583
$ws.Range("M5").Value = -0.1980138146036305
$ws.Range("N5").Value = 0.4101233346279459
$ws.Range("O5").Value = 0.09082176799679093

# Row 7
$ws.Range("B7").Value = -0.02527535427880465
$ws.Range("C7").Value = -0.06051795811232875
$ws.Range("D7").Value = 0.3269328383725711
$ws.Range("E7").Value = 0.952956949937701
$ws.Range("H7").Value = 0.3122748691730542
$ws.Range("I7").Value = 0.3819806403851467
$ws.Range("J7").Value = 0.6010659765620737
$ws.Range("K7").Value = 0.1802842195832685
$ws.Range("L7").Value = -0.0427452430127298
$ws.Range("M7").Value = -0.07627511032712242
$ws.Range("N7").Value = 0.3858841094161588
$ws.Range("O7").Value = 0.1521262490988143

# Row 8
$ws.Range("B8").Value = -0.3848572087141396
$ws.Range("C8").Value = 0.1143498427673656
$ws.Range("D8").Value = 0.3958692142618715
$ws.Range("E8").Value = 0.3616315137833193
$ws.Range("G8").Value = 0.3122748691730542
$ws.Range("I8").Value = -0.2083907738613129
$ws.Range("J8").Value = 0.445985760966529
$ws.Range("K8").Value = -0.1506489769386493
$ws.Range("L8").Value = -0.3230280538554612
$ws.Range("M8").Value = -0.225486848726541
$ws.Range("N8").Value = 0.9082083413765707
$ws.Range("O8").Value = -0.1723749541018063

# Row 9
$ws.Range("B9").Value = 0.07613113137254446
$ws.Range("C9").Value = 0.09527665589281313
$ws.Range("D9").Value = 0.03844707986648056
$ws.Range("E9").Value = 0.3244402953268535
$ws.Range("G9").Value = 0.3819806403851467
$ws.Range("H9").Value = -0.2083907738613129
$ws.Range("J9").Value = 0.1987560697275382
$ws.Range("K9").Value = 0.6522674128189277
$ws.Range("L9").Value = 0.114092003527871
$ws.Range("M9").Value = 0.4944986209532007
$ws.Range("N9").Value = -0.1579486423221279
$ws.Range("O9").Value = 0.7403132706269152

# Row 10
$ws.Range("B10").Value = 0.09692859617640064
$ws.Range("C10").Value = -0.01151755191354742
$ws.Range("D10").Value = 0.8041894500608463
$ws.Range("E10").Value = 0.6315695168996384
$ws.Range("G10").Value = 0.6010659765620737
$ws.Range("H10").Value = 0.445985760966529
$ws.Range("I10").Value = 0.1987560697275382
$ws.Range("K10").Value = 0.2525430715716202
$ws.Range("L10").Value = -0.1349869902559086
$ws.Range("M10").Value = 0.05601955201777537
$ws.Range("N10").Value = 0.5188514403663723
$ws.Range("O10").Value = 0.1893413854013475

# Row 11
$ws.Range("B11").Value = 0.1457023439347263
$ws.Range("C11").Value = -0.02914244125909981
$ws.Range("D11").Value = 0.1418704476126033
$ws.Range("E11").Value = 0.06432798670309411
$ws.Range("G11").Value = 0.1802842195832685
$ws.Range("H11").Value = -0.1506489769386493
$ws.Range("I11").Value = 0.6522674128189277
$ws.Range("J11").Value = 0.2525430715716202
$ws.Range("L11").Value = 0.1177035147288679
$ws.Range("M11").Value = 0.9363545332786369
$ws.Range("N11").Value = -0.07480865452710381
$ws.Range("O11").Value = 0.7459856900792554

# Row 12
$ws.Range("B12").Value = 0.06131781591624236
$ws.Range("C12").Value = 0.04876962339280161
$ws.Range("D12").Value = -0.01463024890903251
$ws.Range("E12").Value = -0.09781173102871583
$ws.Range("G12").Value = -0.0427452430127298
$ws.Range("H12").Value = -0.3230280538554612
$ws.Range("I12").Value = 0.114092003527871
$ws.Range("J12").Value = -0.1349869902559086
$ws.Range("K12").Value = 0.1177035147288679
$ws.Range("M12").Value = 0.1285181598510687
$ws.Range("N12").Value = -0.2839976291441505
$ws.Range("O12").Value = -0.02349543624334667

# Row 13
$ws.Range("B13").Value = 0.06612057560974671
$ws.Range("C13").Value = 0.01742559811962423
$ws.Range("D13").Value = 0.05679644052597838
$ws.Range("E13").Value = -0.1980138146036305
$ws.Range("G13").Value = -0.07627511032712242
$ws.Range("H13").Value = -0.225486848726541
$ws.Range("I13").Value = 0.4944986209532007
$ws.Range("J13").Value = 0.05601955201777537
$ws.Range("K13").Value = 0.9363545332786369
$ws.Range("L13").Value = 0.1285181598510687
$ws.Range("N13").Value = -0.1775736450325868
$ws.Range("O13").Value = 0.6916111360249111

# Row 14
$ws.Range("B14").Value = -0.3919766247436938
$ws.Range("C14").Value = 0.06934592894886464
$ws.Range("D14").Value = 0.3661331865400339
$ws.Range("E14").Value = 0.4101233346279459
$ws.Range("G14").Value = 0.3858841094161588
$ws.Range("H14").Value = 0.9082083413765707
$ws.Range("I14").Value = -0.1579486423221279
$ws.Range("J14").Value = 0.5188514403663723
$ws.Range("K14").Value = -0.07480865452710381
$ws.Range("L14").Value = -0.2839976291441505
$ws.Range("M14").Value = -0.1775736450325868
$ws.Range("O14").Value = -0.1765699855808679

# Row 15
$ws.Range("B15").Value = 0.1208259260153849
$ws.Range("C15").Value = 0.06306186227928907
$ws.Range("D15").Value = 0.04245051279411372
$ws.Range("E15").Value = 0.09082176799679093
$ws.Range("G15").Value = 0.1521262490988143
$ws.Range("H15").Value = -0.1723749541018063
$ws.Range("I15").Value = 0.7403132706269152
$ws.Range("J15").Value = 0.1893413854013475
$ws.Range("K15").Value = 0.7459856900792554
$ws.Range("L15").Value = -0.02349543624334667
$ws.Range("M15").Value = 0.6916111360249111
$ws.Range("N15").Value = -0.1765699855808679
